# TC02_Canine_Filter_Study-NCATS.xlsx
# Insert a new "StatQuery" column between the existing "query" (A) and
# "dbExcel" (B) columns, shifting dbExcel/WebExcel one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B (old B "dbExcel" and C "WebExcel" shift right)
$ws.Columns("B:B").Insert()

# New column should be as wide as column A
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Header for the new column
$ws.Range("B1").Value = "StatQuery"

# New query text (matches the long-form count/stat query added to sharedStrings)
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN ['NCATS-COP01']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Row 2 wraps text (same style as A2) for the new cell
$ws.Range("B2").WrapText = $true

# Restore selection on A2 (single cell, not full-column)
[void]$ws.Range("A2").Select()
